# Add 2022-Q3 data:
#  - insert a new "2022-Q3" worksheet right after "总计", cloned (with its
#    formatting) from the existing "2022-Q1" sheet, then overwrite with the
#    new quarter's fund-holding rows (only 2 rows of data this quarter).
#  - add a summary row for 2022-Q3 at the top of the "总计" sheet's data,
#    pushing the existing quarter rows down by one.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)        # "总计"
$q1Sheet = $wb.Worksheets.Item(2)       # "2022-Q1" (existing, keeps its data)

# --- 1. Create the new "2022-Q3" sheet right after "总计" ------------------
# Cloning "2022-Q1" preserves the header row + styles (bold/centered/bordered
# index & header cells) so the new sheet matches the look of its siblings.
$q1Sheet.Copy($null, $total)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# This quarter only has 2 holdings, so drop the cloned 3rd data row (row 4).
$q3Sheet.Rows.Item(4).Delete()

# Row 2: 164811 / 工银瑞信中证京津冀协同发展主题指数（LOF）A
$q3Sheet.Range("B2").NumberFormat = "@"
$q3Sheet.Range("B2").Value = "164811"
$q3Sheet.Range("C2").Value = "工银瑞信中证京津冀协同发展主题指数（LOF）A"
$q3Sheet.Range("D2").NumberFormat = "@"
$q3Sheet.Range("D2").Value = "0.12"
$q3Sheet.Range("E2").NumberFormat = "@"
$q3Sheet.Range("E2").Value = "93.09"
$q3Sheet.Range("F2").NumberFormat = "@"
$q3Sheet.Range("F2").Value = "2.99"
$q3Sheet.Range("G2").NumberFormat = "@"
$q3Sheet.Range("G2").Value = "0.0036"
$q3Sheet.Range("H2").Value = 7

# Row 3: 164825 / 工银瑞信中证京津冀协同发展主题指数（LOF）C
$q3Sheet.Range("B3").NumberFormat = "@"
$q3Sheet.Range("B3").Value = "164825"
$q3Sheet.Range("C3").Value = "工银瑞信中证京津冀协同发展主题指数（LOF）C"
$q3Sheet.Range("D3").NumberFormat = "@"
$q3Sheet.Range("D3").Value = "0.03"
$q3Sheet.Range("E3").NumberFormat = "@"
$q3Sheet.Range("E3").Value = "93.09"
$q3Sheet.Range("F3").NumberFormat = "@"
$q3Sheet.Range("F3").Value = "2.99"
$q3Sheet.Range("G3").NumberFormat = "@"
$q3Sheet.Range("G3").Value = "0.0009"
$q3Sheet.Range("H3").Value = 7

# --- 2. Insert the new summary row on "总计" --------------------------------
$total.Rows.Item(2).Insert()
# Insert() carries the header row's formatting onto the new row; strip it
# from B:D (those columns hold plain data, no style, on every other row)
# and instead clone column A's index-cell style (bold/centered/bordered).
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0
